$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Insert Phu's new "Software" slide as slide 7 (before the existing
#    "Schedule" slide, which shifts down to position 8). Use the same
#    "Title and Content" layout as the surrounding slides.
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add(7, 2)

$newTitle = $newSlide.Shapes.Item(1).TextFrame.TextRange
$newTitle.Text = "Software"

$newBody = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newBody.Text = "System:`rControl the operation:`rAs per programed schedule -Interrupt from RTC`rManage the data storage `rStore data to SD card`rReset power when SD card " + [char]8220 + "full" + [char]8221 + " for switching to the next available SD card`rManage the interfere between programmed schedule (RTC interrupt) and storage process.`r`r"

$newBody.Paragraphs(2).IndentLevel = 2
$newBody.Paragraphs(3).IndentLevel = 3
$newBody.Paragraphs(4).IndentLevel = 3
$newBody.Paragraphs(5).IndentLevel = 4
$newBody.Paragraphs(6).IndentLevel = 4
$newBody.Paragraphs(7).IndentLevel = 3
$newBody.Paragraphs(8).IndentLevel = 4
$newBody.Paragraphs(9).IndentLevel = 2

# ---------------------------------------------------------------------
# 2. Hardware slide (slide 5): fill in the previously empty content
#    placeholder with the software/UI requirements bullets.
# ---------------------------------------------------------------------
$hw = $p.Slides.Item(5).Shapes.Item(2).TextFrame.TextRange
$hw.Text = "User interface - command line, better " + [char]8211 + " GUI`rUser" + [char]8217 + "s input : recording and stand-by interval, sample rate, delay for start`rResponse: `rSet-up configuration if the inputs are " + [char]8220 + "good" + [char]8221 + "`rEstimated maximum operation time of the system`rAsk for confirmation`rOutput to the system the confirmed operation set-points`rDisplay the signal spectogram`r"

$hw.Paragraphs(2).IndentLevel = 2
$hw.Paragraphs(3).IndentLevel = 2
$hw.Paragraphs(4).IndentLevel = 3
$hw.Paragraphs(5).IndentLevel = 3
$hw.Paragraphs(6).IndentLevel = 3
$hw.Paragraphs(7).IndentLevel = 2
$hw.Paragraphs(8).IndentLevel = 2
$hw.Paragraphs(9).IndentLevel = 2

# ---------------------------------------------------------------------
# 3. Feature Goals (Cont.) slide (slide 4): merge the split
#    "Must create battery " / "harness" runs into a single run.
# ---------------------------------------------------------------------
$fg = $p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange
$fg.Paragraphs(4).Text = "temp-merge-placeholder"
$fg.Paragraphs(4).Text = "Must create battery harness"

Write-Output "done"
